# Split the TEST48 description run into three runs so that a comma is
# inserted right after "test28" and the existing "_GoBack" bookmark ends up
# sitting between the new comma run and the " but using the heat solver. "
# tail run, e.g.:
#
#   " Solves same problem as test28" + "," + <bookmark> + " but using the heat solver. "
#
$d = $word.ActiveDocument

# --- locate the target paragraph --------------------------------------
$marker = "Solves same problem as test28"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.Contains($marker)) {
        $target = $pp
        break
    }
}

$pStart = $target.Range.Start
$paraText = $target.Range.Text

# Offset (relative to paragraph start) of the character right after "test28"
$splitWord = "test28"
$relIdx = $paraText.IndexOf($splitWord)
$absAfterTest28 = $pStart + $relIdx + $splitWord.Length

# --- remove the old "_GoBack" bookmark (it will be re-added later) ----
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- drop everything after "test28" (i.e. " but using the heat solver. ")
$rTail = $d.Range($absAfterTest28, $target.Range.End)
$rTail.Text = ""

# --- append the comma as its own run -----------------------------------
# Inserting right at the paragraph-end boundary (rather than mid-run)
# keeps it from being merged back into the preceding run.
$pEnd = $target.Range.End
$rComma = $d.Range($pEnd - 1, $pEnd - 1)
$rComma.InsertAfter(",")

# --- append the tail text as its own run too, still at the boundary ---
$pEnd2 = $target.Range.End
$rTailNew = $d.Range($pEnd2 - 1, $pEnd2 - 1)
$rTailNew.InsertAfter(" but using the heat solver. ")

# --- re-insert the "_GoBack" bookmark between the comma and the tail --
# (the boundary sits right after the comma run, which by now is a normal
# mid-paragraph position rather than the paragraph's last slot, so
# Bookmarks.Add behaves correctly here)
$bmPos = $pEnd2 - 1
$rBm = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $rBm)
